{"js": "// The use-case table's \"\u57fa\u672c\u7cfb\u5217\" (basic flow) step currently reads:\n//   \"\u30a2\u30af\u30bf\u30fc\u306f\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\u3092\u5165\u529b\u3057\u300c\u78ba\u8a8d\u753b\u9762\u3078\u300d\u30dc\u30bf\u30f3\u3092\u62bc\u3059\"\n// It needs to become:\n//   \"\u30a2\u30af\u30bf\u30fc\u306f\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\u306e\u540d\u524d\u3001\u4f4f\u6240\u3001\u96fb\u8a71\u756a\u53f7\u3001\u30e1\u30fc\u30eb\u30a2\u30c9\u30ec\u30b9\u3001\u751f\u5e74\u6708\u65e5\u306e\u3044\u305a\u308c\u304b\u3092\n//    \u5165\u529b\u3057\u300c\u78ba\u8a8d\u753b\u9762\u3078\u300d\u30dc\u30bf\u30f3\u3092\u62bc\u3059\"\n// i.e. insert \"\u306e\" + \"\u540d\u524d\u3001\u4f4f\u6240\u3001\u96fb\u8a71\u756a\u53f7\u3001\u30e1\u30fc\u30eb\u30a2\u30c9\u30ec\u30b9\u3001\u751f\u5e74\u6708\u65e5\" + \"\u306e\u3044\u305a\u308c\u304b\"\n// immediately after the existing \"\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\" text, right before \"\u3092\u5165\u529b\u3057\".\n\nconst results = context.document.body.search(\"\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\" in the document body.');\n}\n\n// Anchor on the (unique) match and append the three new pieces of text\n// immediately after it, chaining off each inserted range so the pieces land\n// in the same reading order as the source diff.\nlet insertionPoint = results.items[0];\ninsertionPoint = insertionPoint.insertText(\"\u306e\", Word.InsertLocation.after);\ninsertionPoint = insertionPoint.insertText(\n  \"\u540d\u524d\u3001\u4f4f\u6240\u3001\u96fb\u8a71\u756a\u53f7\u3001\u30e1\u30fc\u30eb\u30a2\u30c9\u30ec\u30b9\u3001\u751f\u5e74\u6708\u65e5\",\n  Word.InsertLocation.after\n);\ninsertionPoint = insertionPoint.insertText(\"\u306e\u3044\u305a\u308c\u304b\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# The use-case table's \"\u57fa\u672c\u7cfb\u5217\" (basic flow) step currently reads:\n#   \"\u30a2\u30af\u30bf\u30fc\u306f\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\u3092\u5165\u529b\u3057\u300c\u78ba\u8a8d\u753b\u9762\u3078\u300d\u30dc\u30bf\u30f3\u3092\u62bc\u3059\"\n# It needs to become:\n#   \"\u30a2\u30af\u30bf\u30fc\u306f\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\u306e\u540d\u524d\u3001\u4f4f\u6240\u3001\u96fb\u8a71\u756a\u53f7\u3001\u30e1\u30fc\u30eb\u30a2\u30c9\u30ec\u30b9\u3001\u751f\u5e74\u6708\u65e5\u306e\u3044\u305a\u308c\u304b\u3092\n#    \u5165\u529b\u3057\u300c\u78ba\u8a8d\u753b\u9762\u3078\u300d\u30dc\u30bf\u30f3\u3092\u62bc\u3059\"\n# i.e. insert \"\u306e\" + \"\u540d\u524d\u3001\u4f4f\u6240\u3001\u96fb\u8a71\u756a\u53f7\u3001\u30e1\u30fc\u30eb\u30a2\u30c9\u30ec\u30b9\u3001\u751f\u5e74\u6708\u65e5\" + \"\u306e\u3044\u305a\u308c\u304b\"\n# immediately after the existing \"\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\" text, right before \"\u3092\u5165\u529b\u3057\".\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\")\nif (-not $found) {\n    throw 'Could not find \"\u691c\u7d22\u3057\u305f\u3044\u4f1a\u54e1\" in the document.'\n}\n\n# Collapse to the end of the match and insert the three new pieces of text\n# right after it, advancing the insertion range each time so they appear in\n# order, immediately before \"\u3092\u5165\u529b\u3057\".\n$rng.Collapse(0)\n$rng.InsertAfter(\"\u306e\")\n$rng.Collapse(0)\n$rng.InsertAfter(\"\u540d\u524d\u3001\u4f4f\u6240\u3001\u96fb\u8a71\u756a\u53f7\u3001\u30e1\u30fc\u30eb\u30a2\u30c9\u30ec\u30b9\u3001\u751f\u5e74\u6708\u65e5\")\n$rng.Collapse(0)\n$rng.InsertAfter(\"\u306e\u3044\u305a\u308c\u304b\")\n"}
